# error solve ifrs list
# Rewrites the per-period financial figures on the "company_list" sheet.
# Rows 2-6 (FY2014-FY2018 IFRS-consolidated) get corrected figures, and the
# "자산총계"(J)/"자본총계(비지배)"(O) columns - which should not have been
# populated for this ticker - are cleared. Rows 7-9 (the 2019E/2020E/2021E
# estimate rows) had bad data pulled in and are cleared back to just their
# period label (columns D..AI removed) pending a good refill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12 IFRS연결) ---
$ws.Range("D2").Value = 3274
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = -1
$ws.Range("H2").Value = -89
$ws.Range("I2").Value = -89
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 3496
$ws.Range("L2").Value = 2207
$ws.Range("M2").Value = 1289
$ws.Range("N2").Value = 1289
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 327
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = -186
$ws.Range("S2").Value = 189
$ws.Range("T2").Value = 129
$ws.Range("U2").Value = -117
$ws.Range("V2").Value = 1355
$ws.Range("W2").Value = 0.4
$ws.Range("X2").Value = -2.72
$ws.Range("Y2").Value = -6.61
$ws.Range("Z2").Value = -2.59
$ws.Range("AA2").Value = 171.31
$ws.Range("AB2").Value = 296.45
$ws.Range("AC2").Value = -113
$ws.Range("AD2").Value = -5.58
$ws.Range("AE2").Value = 1632
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 8
$ws.Range("AH2").Value = 1.32
$ws.Range("AI2").Value = -5.53
$ws.Range("AJ2").Value = 78956148

# --- Row 3 (2015/12 IFRS연결) ---
$ws.Range("D3").Value = 3575
$ws.Range("E3").Value = 53
$ws.Range("F3").Value = 67
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 27
$ws.Range("I3").Value = 27
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 3418
$ws.Range("L3").Value = 2117
$ws.Range("M3").Value = 1302
$ws.Range("N3").Value = 1302
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 327
$ws.Range("Q3").Value = 184
$ws.Range("R3").Value = -124
$ws.Range("S3").Value = -62
$ws.Range("T3").Value = 58
$ws.Range("U3").Value = 125
$ws.Range("V3").Value = 1298
$ws.Range("W3").Value = 1.48
$ws.Range("X3").Value = 0.76
$ws.Range("Y3").Value = 2.09
$ws.Range("Z3").Value = 0.78
$ws.Range("AA3").Value = 162.61
$ws.Range("AB3").Value = 298.18
$ws.Range("AC3").Value = 34
$ws.Range("AD3").Value = 29.96
$ws.Range("AE3").Value = 1649
$ws.Range("AF3").Value = 0.62
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 1.21
$ws.Range("AI3").Value = 36.27
$ws.Range("AJ3").Value = 78956148

# --- Row 4 (2016/12 IFRS연결) ---
$ws.Range("D4").Value = 3589
$ws.Range("E4").Value = -11
$ws.Range("F4").Value = -11
$ws.Range("G4").Value = -4
$ws.Range("H4").Value = 89
$ws.Range("I4").Value = 28
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 5620
$ws.Range("L4").Value = 4290
$ws.Range("M4").Value = 1329
$ws.Range("N4").Value = 844
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 327
$ws.Range("Q4").Value = 114
$ws.Range("R4").Value = -720
$ws.Range("S4").Value = 538
$ws.Range("T4").Value = 115
$ws.Range("U4").Value = -2
$ws.Range("V4").Value = 2150
$ws.Range("W4").Value = -0.32
$ws.Range("X4").Value = 2.49
$ws.Range("Y4").Value = 2.64
$ws.Range("Z4").Value = 1.98
$ws.Range("AA4").Value = 322.71
$ws.Range("AB4").Value = 290.14
$ws.Range("AC4").Value = 36
$ws.Range("AD4").Value = 87.58
$ws.Range("AE4").Value = 1069
$ws.Range("AF4").Value = 2.93
$ws.Range("AG4").Value = 8
$ws.Range("AH4").Value = 0.26
$ws.Range("AI4").Value = 23.12
$ws.Range("AJ4").Value = 78956148

# --- Row 5 (2017/12 IFRS연결) ---
$ws.Range("D5").Value = 5283
$ws.Range("E5").Value = -19
$ws.Range("F5").Value = -19
$ws.Range("G5").Value = -239
$ws.Range("H5").Value = -226
$ws.Range("I5").Value = -223
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 5366
$ws.Range("L5").Value = 4279
$ws.Range("M5").Value = 1088
$ws.Range("N5").Value = 610
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 327
$ws.Range("Q5").Value = 212
$ws.Range("R5").Value = -133
$ws.Range("S5").Value = -82
$ws.Range("T5").Value = 108
$ws.Range("U5").Value = 104
$ws.Range("V5").Value = 2061
$ws.Range("W5").Value = -0.36
$ws.Range("X5").Value = -4.28
$ws.Range("Y5").Value = -30.72
$ws.Range("Z5").Value = -4.12
$ws.Range("AA5").Value = 393.38
$ws.Range("AB5").Value = 85.61
$ws.Range("AC5").Value = -283
$ws.Range("AD5").Value = -4.77
$ws.Range("AE5").Value = 773
$ws.Range("AF5").Value = 1.75
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 0.92
$ws.Range("AI5").Value = -4.39
$ws.Range("AJ5").Value = 78956148

# --- Row 6 (2018/12 IFRS연결) ---
# (this row never had J/O populated, so there is nothing to clear there)
$ws.Range("D6").Value = 5899
$ws.Range("E6").Value = -54
$ws.Range("F6").Value = -54
$ws.Range("G6").Value = -138
$ws.Range("H6").Value = -228
$ws.Range("I6").Value = -75
$ws.Range("K6").Value = 5815
$ws.Range("L6").Value = 4773
$ws.Range("M6").Value = 1043
$ws.Range("N6").Value = 665
$ws.Range("P6").Value = 327
$ws.Range("Q6").Value = 207
$ws.Range("R6").Value = -107
$ws.Range("S6").Value = -101
$ws.Range("T6").Value = 139
$ws.Range("U6").Value = 68
$ws.Range("V6").Value = 1969
$ws.Range("W6").Value = -0.92
$ws.Range("X6").Value = -3.86
$ws.Range("Y6").Value = -11.73
$ws.Range("Z6").Value = -4.07
$ws.Range("AA6").Value = 457.77
$ws.Range("AB6").Value = 53.06
$ws.Range("AC6").Value = -95
$ws.Range("AD6").Value = -15.94
$ws.Range("AE6").Value = 842
$ws.Range("AF6").Value = 1.79
$ws.Range("AG6").Value = 8
$ws.Range("AH6").Value = 0.55
$ws.Range("AI6").Value = -6.57
$ws.Range("AJ6").Value = 78956148

# --- Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# Bad estimate data - blank out all the figure columns, leaving only the
# row index / ticker / period label columns (A:C) intact.
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
